$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.230798482894897
$ws.Range("B1").Value = 2.52327823638916
$ws.Range("C1").Value = 9.187633514404297
$ws.Range("D1").Value = 2.043156385421753
$ws.Range("E1").Value = 1.1752769947052
